# Auto-generated Excel COM-interop script to apply the Kraken_Profits diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 248.45454
$ws.Range("I9").Value = 258.8
$ws.Range("K9").Value = 258.8
$ws.Range("M9").Value = -89.80000000000001
# Row 29
$ws.Range("H29").Value = 450
$ws.Range("I29").Value = 187.5
$ws.Range("K29").Value = 562.5
$ws.Range("M29").Value = -281.5
# Row 33
$ws.Range("H33").Value = 288.66666
$ws.Range("I33").Value = 288.66666
$ws.Range("K33").Value = 288.66666
$ws.Range("M33").Value = -59.66665999999998
# Row 70
$ws.Range("H70").Value = 2858.4167
$ws.Range("I70").Value = 1900
$ws.Range("J70").Value = 2945.5454
$ws.Range("K70").Value = 5700
$ws.Range("L70").Value = 8836.636200000001
$ws.Range("M70").Value = -5430
$ws.Range("N70").Value = -9376.636200000001
# Row 73
$ws.Range("H73").Value = 2858.4167
$ws.Range("I73").Value = 1900
$ws.Range("J73").Value = 2945.5454
$ws.Range("K73").Value = 5700
$ws.Range("L73").Value = 8836.636200000001
$ws.Range("M73").Value = -4764
$ws.Range("N73").Value = -10708.6362
# Row 80
$ws.Range("H80").Value = 349.83334
$ws.Range("I80").Value = 299.75
$ws.Range("J80").Value = 450
$ws.Range("K80").Value = 899.25
$ws.Range("L80").Value = 1350
$ws.Range("M80").Value = 98.75
$ws.Range("N80").Value = -3346
# Row 83
$ws.Range("H83").Value = 349.83334
$ws.Range("I83").Value = 299.75
$ws.Range("J83").Value = 450
$ws.Range("K83").Value = 2697.75
$ws.Range("L83").Value = 4050
$ws.Range("M83").Value = 2294.25
$ws.Range("N83").Value = -14034
# Row 125
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6540
$ws.Range("N125").ClearContents()
# Row 137
$ws.Range("H137").Value = 2247
$ws.Range("I137").Value = 2163
$ws.Range("K137").Value = 6489
$ws.Range("M137").Value = -3939
# Row 138
$ws.Range("H138").Value = 3619.5
$ws.Range("J138").Value = 3998.95
$ws.Range("L138").Value = 11996.85
$ws.Range("N138").Value = -22276.85

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 43
$ws.Range("H43").Value = 47823.5
$ws.Range("I43").Value = 47824
$ws.Range("J43").Value = 47823
$ws.Range("K43").Value = 47824
$ws.Range("L43").Value = 47823
$ws.Range("M43").Value = -47511
$ws.Range("N43").Value = -48449
# Row 74
$ws.Range("H74").Value = 2342.4546
$ws.Range("I74").Value = 2199.7778
$ws.Range("J74").Value = 2984.5
$ws.Range("K74").Value = 2199.7778
$ws.Range("L74").Value = 2984.5
$ws.Range("M74").Value = -1325.7778
$ws.Range("N74").Value = -4732.5
# Row 77
$ws.Range("H77").Value = 2342.4546
$ws.Range("I77").Value = 2199.7778
$ws.Range("J77").Value = 2984.5
$ws.Range("K77").Value = 10998.889
$ws.Range("L77").Value = 14922.5
$ws.Range("M77").Value = -6630.888999999999
$ws.Range("N77").Value = -23658.5
# Row 110
$ws.Range("H110").Value = 3133.3333
$ws.Range("J110").Value = 7000
$ws.Range("L110").Value = 7000
$ws.Range("N110").Value = -11090
# Row 112
$ws.Range("H112").Value = 38877.2
$ws.Range("J112").Value = 38877.2
$ws.Range("L112").Value = 38877.2
$ws.Range("N112").Value = -41831.2
# Row 122
$ws.Range("H122").Value = 6865.759
$ws.Range("I122").Value = 6264.64
$ws.Range("K122").Value = 18793.92
$ws.Range("M122").Value = -16343.92
# Row 124
$ws.Range("H124").Value = 42500
$ws.Range("J124").Value = 42500
$ws.Range("L124").Value = 42500
$ws.Range("N124").Value = -52320
# Row 132
$ws.Range("H132").Value = 1879.8572
$ws.Range("I132").Value = 1779.3684
$ws.Range("J132").Value = 2834.5
$ws.Range("K132").Value = 5338.1052
$ws.Range("L132").Value = 8503.5
$ws.Range("M132").Value = -2808.1052
$ws.Range("N132").Value = -13563.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2362.125
$ws.Range("I86").Value = 2582.8333
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 2582.8333
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -1459.8333
$ws.Range("N86").Value = -3946
# Row 89
$ws.Range("H89").Value = 2362.125
$ws.Range("I89").Value = 2582.8333
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 12914.1665
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -7298.166499999999
$ws.Range("N89").Value = -19732
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
# Row 99
$ws.Range("H99").Value = 5272.091
$ws.Range("I99").Value = 5272.091
$ws.Range("K99").Value = 5272.091
$ws.Range("M99").Value = -3774.091
# Row 105
$ws.Range("H105").Value = 4993
$ws.Range("I105").Value = 4993
$ws.Range("K105").Value = 4993
$ws.Range("M105").Value = -3246
# Row 134
$ws.Range("H134").Value = 5379.5293
$ws.Range("I134").Value = 1787
$ws.Range("J134").Value = 14001.6
$ws.Range("K134").Value = 5361
$ws.Range("L134").Value = 42004.8
$ws.Range("M134").Value = -2826
$ws.Range("N134").Value = -47074.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 33.266666
$ws.Range("I12").Value = 100.5
$ws.Range("J12").Value = 22.923077
$ws.Range("K12").Value = 301.5
$ws.Range("L12").Value = 68.76923099999999
$ws.Range("M12").Value = -128.5
$ws.Range("N12").Value = -414.769231
# Row 19
$ws.Range("H19").Value = 3999
$ws.Range("J19").Value = 3999
$ws.Range("L19").Value = 11997
$ws.Range("N19").Value = -12345
# Row 108
$ws.Range("H108").Value = 760
$ws.Range("I108").Value = 160
$ws.Range("K108").Value = 480
$ws.Range("M108").Value = 2400
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# Row 111
$ws.Range("H111").Value = 1998.5
$ws.Range("I111").Value = 1998.5
$ws.Range("K111").Value = 5995.5
$ws.Range("M111").Value = -2928.5
# Row 112
$ws.Range("H112").Value = 999
$ws.Range("I112").Value = 999
$ws.Range("K112").Value = 2997
$ws.Range("M112").Value = -1889
# Row 128
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -10996
# Row 83
$ws.Range("H83").Value = 9000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 45000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -54984
# Row 97
$ws.Range("H97").Value = 611.3570999999999
$ws.Range("I97").Value = 617.3333
$ws.Range("K97").Value = 617.3333
$ws.Range("M97").Value = -121.3333
# Row 107
$ws.Range("H107").Value = 659.9091
$ws.Range("I107").Value = 562.2222
$ws.Range("K107").Value = 562.2222
$ws.Range("M107").Value = 1357.7778
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
# Row 122
$ws.Range("H122").Value = 2832.2
$ws.Range("I122").Value = 3073.7778
$ws.Range("J122").Value = 658
$ws.Range("K122").Value = 9221.3334
$ws.Range("L122").Value = 1974
$ws.Range("M122").Value = -6771.3334
$ws.Range("N122").Value = -6874
# Row 132
$ws.Range("H132").Value = 5097.1
$ws.Range("I132").Value = 4495.8335
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 13487.5005
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -10957.5005
$ws.Range("N132").Value = -23057

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3460
$ws.Range("I61").Value = 2825
$ws.Range("J61").Value = 6000
$ws.Range("K61").Value = 2825
$ws.Range("L61").Value = 6000
$ws.Range("M61").Value = -2623
$ws.Range("N61").Value = -6404
# Row 113
$ws.Range("H113").Value = 3460
$ws.Range("I113").Value = 2825
$ws.Range("J113").Value = 6000
$ws.Range("K113").Value = 2825
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -655
$ws.Range("N113").Value = -10340

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5159.25
$ws.Range("I62").Value = 7333.3335
$ws.Range("J62").Value = 3854.8
$ws.Range("K62").Value = 7333.3335
$ws.Range("L62").Value = 3854.8
$ws.Range("M62").Value = -6709.3335
$ws.Range("N62").Value = -5102.8
# Row 65
$ws.Range("H65").Value = 5159.25
$ws.Range("I65").Value = 7333.3335
$ws.Range("J65").Value = 3854.8
$ws.Range("K65").Value = 36666.6675
$ws.Range("L65").Value = 19274
$ws.Range("M65").Value = -33546.6675
$ws.Range("N65").Value = -25514
# Row 81
$ws.Range("H81").Value = 1575.4
$ws.Range("I81").Value = 1575.4
$ws.Range("K81").Value = 3150.8
$ws.Range("M81").Value = -2089.8
# Row 84
$ws.Range("H84").Value = 1575.4
$ws.Range("I84").Value = 1575.4
$ws.Range("K84").Value = 15754
$ws.Range("M84").Value = -10450
# Row 107
$ws.Range("H107").Value = 1289.6
$ws.Range("I107").Value = 274.5
$ws.Range("J107").Value = 1966.3334
$ws.Range("K107").Value = 823.5
$ws.Range("L107").Value = 5899.0002
$ws.Range("M107").Value = 1096.5
$ws.Range("N107").Value = -9739.0002
# Row 112
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
# Row 122
$ws.Range("H122").Value = 2431.1765
$ws.Range("I122").Value = 2023.6428
$ws.Range("J122").Value = 4333
$ws.Range("K122").Value = 6070.928400000001
$ws.Range("L122").Value = 12999
$ws.Range("M122").Value = -3620.928400000001
$ws.Range("N122").Value = -17899
# Row 126
$ws.Range("H126").Value = 3058
$ws.Range("I126").Value = 2322.5
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 6967.5
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -4497.5
$ws.Range("N126").Value = -22940
# Row 132
$ws.Range("H132").Value = 2964.2856
$ws.Range("I132").Value = 1461.9412
$ws.Range("J132").Value = 9349.25
$ws.Range("K132").Value = 4385.8236
$ws.Range("L132").Value = 28047.75
$ws.Range("M132").Value = -1855.8236
$ws.Range("N132").Value = -33107.75

